$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

# Check the two form-control checkboxes whose linked cells are J26 and J27
# (task rows 26 "Rotating Wall or Door" and 27 "New Look and Feel" -
#  camera following + rotating), and keep their linked cells/value in sync.
$ws.Shapes.Item("Check Box 41").ControlFormat.Value = 1
$ws.Shapes.Item("Check Box 42").ControlFormat.Value = 1

$ws.Range("J26").Value = $true
$ws.Range("J27").Value = $true
